# Daily attendance processing - 2025-12-02 07:29:16
# Normalize the "Recorded By" (column G) entries so that "System" is listed
# before the recorder's email address, e.g.:
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

# Find the last used row in column G ("Recorded By") and sweep every row,
# swapping the value wherever it matches exactly.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

$changed = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
        $changed++
    }
}

Write-Host "Updated $changed 'Recorded By' cell(s) in column G"
